$d = $word.ActiveDocument

$replacements = @(
    @("86÷4=", "44÷7="),
    @("55÷8=", "98÷7="),
    @("87÷9=", "28÷8="),
    @("55÷6=", "73÷3="),
    @("53÷7=", "94÷4="),
    @("47÷9=", "10÷5="),
    @("95÷7=", "37÷3="),
    @("22÷4=", "22÷5="),
    @("73÷2=", "44÷9="),
    @("65÷5=", "88÷9="),
    @("52÷7=", "90÷2="),
    @("17÷7=", "24÷3="),
    @("85÷8=", "88÷6="),
    @("80÷2=", "48÷7="),
    @("40÷4=", "89÷7="),
    @("29÷2=", "43÷6="),
    @("15÷4=", "87÷4="),
    @("98÷3=", "14÷7="),
    @("24÷7=", "15÷3="),
    @("12÷8=", "57÷2="),
    @("11÷3=", "74÷2="),
    @("77÷9=", "22÷5="),
    @("31÷2=", "32÷9="),
    @("36÷7=", "97÷9="),
    @("83÷7=", "13÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
